# Auto-generated Excel COM-interop script to apply betting odds updates
# to Jogos_do_Dia_Betfair_Back_Lay_2026-01-15.xlsx (rows 2-12, columns F..AO)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.71
$ws.Range("G2").Value = 1.74
$ws.Range("H2").Value = 5.7
$ws.Range("J2").Value = 3.85
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 1.45
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 3.5
$ws.Range("O2").Value = 1.37
$ws.Range("P2").Value = 1.84
$ws.Range("Q2").Value = 2.12
$ws.Range("R2").Value = 1.31
$ws.Range("S2").Value = 4
$ws.Range("T2").Value = 2.02
$ws.Range("U2").Value = 1.84
$ws.Range("W2").Value = 2.36
$ws.Range("X2").Value = 13.5
$ws.Range("Y2").Value = 17
$ws.Range("Z2").Value = 46
$ws.Range("AA2").Value = 170
$ws.Range("AB2").Value = 7.4
$ws.Range("AC2").Value = 9.199999999999999
$ws.Range("AD2").Value = 23
$ws.Range("AE2").Value = 100
$ws.Range("AF2").Value = 9.800000000000001
$ws.Range("AI2").Value = 110
$ws.Range("AJ2").Value = 18
$ws.Range("AK2").Value = 21
$ws.Range("AL2").Value = 46
$ws.Range("AM2").Value = 160
$ws.Range("AN2").Value = 13
$ws.Range("AO2").Value = 150

# Row 3
$ws.Range("J3").Value = 8
$ws.Range("K3").Value = 8.4
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 1.1
$ws.Range("P3").Value = 4
$ws.Range("Q3").Value = 1.31
$ws.Range("R3").Value = 2.18
$ws.Range("S3").Value = 1.79
$ws.Range("T3").Value = 1.71
$ws.Range("U3").Value = 2.16
$ws.Range("AB3").Value = 17
$ws.Range("AD3").Value = 46
$ws.Range("AL3").Value = 30
$ws.Range("AN3").Value = 2.96
$ws.Range("AO3").Value = 120

# Row 4
$ws.Range("G4").Value = 800
$ws.Range("J4").Value = 1.04
$ws.Range("N4").Value = 1.26
$ws.Range("P4").Value = 1.26

# Row 5
$ws.Range("G5").Value = 600
$ws.Range("J5").Value = 1.04
$ws.Range("N5").Value = 1.26
$ws.Range("P5").Value = 1.26

# Row 6
$ws.Range("G6").Value = 970
$ws.Range("I6").Value = 970

# Row 7
$ws.Range("F7").Value = 1.04
$ws.Range("H7").Value = 1.04
$ws.Range("J7").Value = 1.09
$ws.Range("N7").Value = 1.1

# Row 8
$ws.Range("G8").Value = 600
$ws.Range("J8").Value = 1.04
$ws.Range("O8").Value = 1.11
$ws.Range("Q8").Value = 1.11
$ws.Range("S8").Value = 1.1

# Row 9
$ws.Range("N9").Value = 2.96
$ws.Range("O9").Value = 1.48
$ws.Range("P9").Value = 1.66
$ws.Range("Q9").Value = 2.42
$ws.Range("R9").Value = 1.24

# Row 10
$ws.Range("N10").Value = 3.15
$ws.Range("O10").Value = 1.45
$ws.Range("U10").Value = 1.98

# Row 11
$ws.Range("S11").Value = 4
$ws.Range("T11").Value = 1.85
$ws.Range("Y11").Value = 9.800000000000001
$ws.Range("AE11").Value = 30
$ws.Range("AO11").Value = 27

# Row 12
$ws.Range("G12").Value = 1.8
$ws.Range("K12").Value = 4.2
$ws.Range("Q12").Value = 1.97
$ws.Range("S12").Value = 3.4
$ws.Range("W12").Value = 2.24
